# Weekly price update: insert the newest week's two rows (Magnum / Sin
# especificar) at the top of the "Poroto verde" data block (row 159),
# pushing the existing historical rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 159:258 down to 161:260, opening up two blank rows.
$ws.Rows("159:160").Insert()

# New row 159 - Magnum, week of 2023-03-16 (serial 45001)
$ws.Cells.Item(159, 1).Value  = 2
$ws.Cells.Item(159, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(159, 3).Value  = "Coquimbo"
$ws.Cells.Item(159, 4).Value  = 45001
$ws.Cells.Item(159, 5).Value  = 4
$ws.Cells.Item(159, 6).Value  = 100112031
$ws.Cells.Item(159, 7).Value  = "Poroto verde"
$ws.Cells.Item(159, 8).Value  = "Magnum"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 800
$ws.Cells.Item(159, 11).Value = 19000
$ws.Cells.Item(159, 12).Value = 20000
$ws.Cells.Item(159, 13).Value = 19500
$ws.Cells.Item(159, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(159, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(159, 16).Value = 780
$ws.Cells.Item(159, 17).Value = 25
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# New row 160 - Sin especificar, same week (serial 45001)
$ws.Cells.Item(160, 1).Value  = 2
$ws.Cells.Item(160, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(160, 3).Value  = "Coquimbo"
$ws.Cells.Item(160, 4).Value  = 45001
$ws.Cells.Item(160, 5).Value  = 4
$ws.Cells.Item(160, 6).Value  = 100112031
$ws.Cells.Item(160, 7).Value  = "Poroto verde"
$ws.Cells.Item(160, 8).Value  = "Sin especificar"
$ws.Cells.Item(160, 9).Value  = "Primera"
$ws.Cells.Item(160, 10).Value = 400
$ws.Cells.Item(160, 11).Value = 23000
$ws.Cells.Item(160, 12).Value = 25000
$ws.Cells.Item(160, 13).Value = 24000
$ws.Cells.Item(160, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 960
$ws.Cells.Item(160, 17).Value = 25
$ws.Cells.Item(160, 18).Value = "Hortaliza"
